$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '50.043.87'
$ws.Range("E2").Value = '  +3.11%  '

# Row 3
$ws.Range("D3").Value = '2.675.36'
$ws.Range("E3").Value = '  +6.24%  '

# Row 4
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.09%  '

# Row 5
$ws.Range("D5").Value = '114.40'
$ws.Range("E5").Value = '  +5.11%  '

# Row 6
$ws.Range("D6").Value = '328.22'
$ws.Range("E6").Value = '  +2.09%  '

# Row 7
$ws.Range("E7").Value = '  +0.90%  '

# Row 8
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.04%  '

# Row 9
$ws.Range("D9").Value = '0.562'
$ws.Range("E9").Value = '  +3.13%  '

# Row 10
$ws.Range("D10").Value = '41.48'
$ws.Range("E10").Value = '  +3.58%  '

# Row 11
$ws.Range("D11").Value = '20.45'
$ws.Range("E11").Value = '  -0.08%  '

# Row 12
$ws.Range("D12").Value = '0.0826'
$ws.Range("E12").Value = '  +1.44%  '

# Row 13
$ws.Range("E13").Value = '  +0.67%  '

# Row 14
$ws.Range("D14").Value = '7.40'
$ws.Range("E14").Value = '  +3.66%  '

# Row 15
$ws.Range("D15").Value = '3.096.06'
$ws.Range("E15").Value = '  +6.30%  '

# Row 16
$ws.Range("D16").Value = '2.680.39'
$ws.Range("E16").Value = '  +6.29%  '

# Row 17
$ws.Range("D17").Value = '0.880'
$ws.Range("E17").Value = '  +4.81%  '

# Row 18
$ws.Range("D18").Value = '50.048.30'
$ws.Range("E18").Value = '  +3.37%  '

# Row 19
$ws.Range("D19").Value = '13.33'
$ws.Range("E19").Value = '  +1.42%  '

# Row 20
$ws.Range("D20").Value = '6.86'
$ws.Range("E20").Value = '  +2.21%  '

# Row 21
$ws.Range("D21").Value = '2.95'
$ws.Range("E21").Value = '  +4.08%  '

# Row 22
$ws.Range("D22").Value = '0.0₃0966'
$ws.Range("E22").Value = '  +2.23%  '

# Row 23
$ws.Range("D23").Value = '281.43'
$ws.Range("E23").Value = '  +0.84%  '

# Row 24
$ws.Range("D24").Value = '73.05'
$ws.Range("E24").Value = '  +1.86%  '

# Row 25
$ws.Range("D25").Value = '2.61'
$ws.Range("E25").Value = '  +2.39%  '

# Row 26
$ws.Range("D26").Value = '27.15'
$ws.Range("E26").Value = '  +4.08%  '

# Row 27
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  -0.07%  '

# Row 28
$ws.Range("B28").Value = 'InjectiveProtocol'
$ws.Range("C28").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D28").Value = '37.06'
$ws.Range("E28").Value = '  +4.45%  '

# Row 29
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").Value = '2.23'
$ws.Range("E29").Value = '  +1.49%  '

# Row 30
$ws.Range("B30").Value = 'Kaspa'
$ws.Range("C30").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D30").Value = '0.144'
$ws.Range("E30").Value = '  -1.56%  '

# Row 31
$ws.Range("B31").Value = 'Cosmos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D31").Value = '9.89'
$ws.Range("E31").Value = '  +1.01%  '

# Row 32
$ws.Range("D32").Value = '50.09'
$ws.Range("E32").Value = '  +0.74%  '

# Row 33
$ws.Range("D33").Value = '19.80'
$ws.Range("E33").Value = '  +1.58%  '

# Row 34
$ws.Range("D34").Value = '5.48'
$ws.Range("E34").Value = '  +2.37%  '

# Row 35
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").Value = '0.0804'
$ws.Range("E35").Value = '  +2.01%  '

# Row 36
$ws.Range("B36").Value = 'FirstDigitalUSD'
$ws.Range("C36").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D36").Value = '1.00'
$ws.Range("E36").Value = '  -0.07%  '

# Row 37
$ws.Range("D37").Value = '2.10'
$ws.Range("E37").Value = '  +7.15%  '

# Row 38
$ws.Range("D38").Value = '4.83'
$ws.Range("E38").Value = '  +3.38%  '

# Row 39
$ws.Range("E39").Value = '  +6.81%  '

# Row 40
$ws.Range("D40").Value = '126.25'
$ws.Range("E40").Value = '  +4.34%  '

# Row 41
$ws.Range("E41").Value = '  +1.30%  '

# Row 42
$ws.Range("B42").Value = 'EnergySwap'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D42").Value = '22.76'
$ws.Range("E42").Value = '  +4.35%  '

# Row 43
$ws.Range("B43").Value = 'WEMIXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D43").Value = '2.26'
$ws.Range("E43").Value = '  +1.70%  '

# Row 44
$ws.Range("D44").Value = '0.0317'
$ws.Range("E44").Value = '  +3.87%  '

# Row 45
$ws.Range("D45").Value = '3.41'
$ws.Range("E45").Value = '  +6.69%  '

# Row 46
$ws.Range("D46").Value = '2.084.09'
$ws.Range("E46").Value = '  +3.23%  '

# Row 47
$ws.Range("E47").Value = '  +12.83%  '

# Row 48
$ws.Range("D48").Value = '2.01'
$ws.Range("E48").Value = '  +6.85%  '

# Row 49
$ws.Range("D49").Value = '9.16'
$ws.Range("E49").Value = '  +1.33%  '

# Row 50
$ws.Range("D50").Value = '5.46'
$ws.Range("E50").Value = '  +3.93%  '

# Row 51
$ws.Range("D51").Value = '83.16'
$ws.Range("E51").Value = '  +3.57%  '
